$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-11 (columns A-I). Empty string "" marks a cell that exists
# but is blank (matches the self-closed <c t="inlineStr"/> elements in the diff).
$data = @(
    @('1000028','Abercrombie, John','Inquiries concerning the intellectual powers and the investigation of truth.','','London :','John Murray,','1838','Location: Cambridge.; Identifier: Abercrombie1838wa66W.; Public number: 0001.',''),
    @('1000029','Darwin, Charles','Inquiries concerning the intellectual powers and the investigation of truth. :',"[Supplementary material in Charles Darwin's copy].",'','','','Identifier: Abercrombie1838wa66W_MS.',''),
    @('1000159','Agassiz, Louis','Lake Superior: its character, vegetation, and animals, compared with those of other similar regions.','','Boston :','Gould, Kendall & Lincoln,','1850','Inscription.; Location: Cambridge.; Identifier: Agassiz1850up52I.; Public number: 0017.',''),
    @('1000160','Darwin, Charles','Lake Superior: its character, vegetation, and animals, compared with those of other similar regions. :',"[Supplementary material in Charles Darwin's copy].",'','','','Identifier: Agassiz1850up52I_MS.',''),
    @('1000223','Agassiz, Louis','Contributions to the natural history of the United States of North America.','','[s.n.] :',',','n.d..','Inscription.; Essay on classification.; Location: Cambridge.; Identifier: Agassiz2006ft69Y.; Public number: 0015.',''),
    @('1000224','Darwin, Charles','Contributions to the natural history of the United States of North America. :',"[Supplementary material in Charles Darwin's copy].",'','','','Identifier: Agassiz2006ft69Y_MS.',''),
    @('1000345','Barker-Webb, Philip','Histoire naturelle des Îles Canaries.','','Paris :','Béthune,','1840','Location: Cambridge.; Identifier: Barker-Webb1840yf41S.; Public number: 0063.',''),
    @('1000346','Darwin, Charles','Histoire naturelle des Îles Canaries. :',"[Supplementary material in Charles Darwin's copy].",'','','','Identifier: Barker-Webb1840yf41S_MS.',''),
    @('1000128','Bechstein, Johann Matthäus','Naturgeschichte der Stubenvögel.','','Halle :','Hennemann,','1840','Signature.; Location: Cambridge.; Identifier: Bechstein1840ob74D.; Public number: 0083.',''),
    @('1000129','Darwin, Charles','Naturgeschichte der Stubenvögel. :',"[Supplementary material in Charles Darwin's copy].",'','','','Identifier: Bechstein1840ob74D_MS.','')
)

$startRow = 2
$endRow = $startRow + $data.Count - 1

# Force the whole block to Text format first, so that:
#  - numeric-looking values (IDs, years) stay stored as text rather than
#    being coerced into numbers, and
#  - cells that end up with an empty string are still materialized in the
#    sheet (as blank cells) instead of being dropped entirely.
$blockRange = $ws.Range("A$startRow`:I$endRow")
$blockRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($rowNum, $col).Value = $rowData[$c]
    }
}

# Restore the default (Normal) style on the block now that all values are
# set, keeping the text-typed values but dropping the explicit style index.
$blockRange.Style = "Normal"

$wb.Save()
